$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.820.36"
$ws.Range("E2").Value = "  +7.54%  "

$ws.Range("D3").Value = "3.637.06"
$ws.Range("E3").Value = "  +7.38%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.92"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.72%  "

$ws.Range("E7").Value = "  +2.84%  "

$ws.Range("D8").Value = "3.618.05"
$ws.Range("E8").Value = "  +7.12%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.80%  "

$ws.Range("E11").Value = "  +4.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.82"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000295"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.75"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.87%  "

$ws.Range("D15").Value = "4.217.09"
$ws.Range("E15").Value = "  +7.46%  "

$ws.Range("D16").Value = "3.637.70"
$ws.Range("E16").Value = "  +7.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.40"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.52%  "

$ws.Range("D18").Value = "70.727.93"
$ws.Range("E18").Value = "  +7.57%  "

$ws.Range("E19").Value = "  +6.40%  "

$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("E21").Value = "  +5.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.55"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.79"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +16.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.47"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +8.85%  "

$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("E27").Value = "  +6.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.25"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.39"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.37"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.72"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +17.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.26"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "614.36"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.76%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.117"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +8.02%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.58"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.56%  "

$ws.Range("D36").Value = "0.0₃0831"
$ws.Range("E36").Value = "  +11.98%  "

$ws.Range("E37").Value = "  +4.68%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.404"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.92%  "

$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.17%  "

$ws.Range("D42").Value = "3.362.64"
$ws.Range("E42").Value = "  +8.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.71%  "

$ws.Range("E44").Value = "  +6.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.66"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.07%  "

$ws.Range("E47").Value = "  +3.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.21"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.06%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.73"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +9.76%  "

$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.25%  "
